$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = 10
$ws.Range("B27").Value = "Vega Modelo de Temuco"
$ws.Range("C27").Value = "La Araucanía"
$ws.Range("D27").Value = 44967
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100107
$ws.Range("H27").Value = "Otros"
$ws.Range("I27").Value = 100107011
$ws.Range("J27").Value = "Tuna"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 80
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 25000
$ws.Range("P27").Value = 25000
$ws.Range("Q27").Value = "$/caja 18 kilos"
$ws.Range("R27").Value = "Provincia de Los Andes"
$ws.Range("S27").Value = 1389
$ws.Range("T27").Value = 18
